$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WSM")

# Row 4: Inventory
$ws.Range("B4").Value = 1006000000.0
$ws.Range("C4").Value = 1125000000.0
$ws.Range("D4").Value = 1042000000.0
$ws.Range("E4").Value = 1071000000.0
$ws.Range("F4").Value = 1101000000.0

# Row 14: Accounts Payable
$ws.Range("B14").Value = 543000000.0
$ws.Range("C14").Value = 562000000.0
$ws.Range("D14").Value = 373000000.0
$ws.Range("E14").Value = 423000000.0
$ws.Range("F14").Value = 521000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -62000000.0
$ws.Range("C22").Value = -42000000.0
$ws.Range("D22").Value = -38000000.0
$ws.Range("E22").Value = -33000000.0
$ws.Range("F22").Value = -48000000.0
